$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.017.55"
$ws.Range("E2").Value = "  +6.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.984.03"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.18"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.51"
$ws.Range("E6").Value = "  +7.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.981.30"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.98"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +3.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.95"
$ws.Range("E14").Value = "  +7.11%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.996.48"
$ws.Range("E16").Value = "  +5.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.475.47"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.91"
$ws.Range("E18").Value = "  +6.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.978.00"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.79"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.68"
$ws.Range("E21").Value = "  +5.07%  "
$ws.Range("E22").Value = "  +4.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.31"
$ws.Range("E23").Value = "  +7.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.99"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.41"
$ws.Range("E25").Value = "  +6.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  +9.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.21"
$ws.Range("E27").Value = "  +10.97%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.43"
$ws.Range("E29").Value = "  +19.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.88"
$ws.Range("E30").Value = "  +13.43%  "
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.77"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.983"
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.72"
$ws.Range("E37").Value = "  +6.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.09"
$ws.Range("E38").Value = "  +8.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "45.30"
$ws.Range("E39").Value = "  +15.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "48.86"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.301"
$ws.Range("E42").Value = "  +13.44%  "
$ws.Range("E43").Value = "  +6.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.40"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "386.43"
$ws.Range("E45").Value = "  +14.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.767.18"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("E47").Value = "  +5.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.20"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.15"
$ws.Range("E50").Value = "  +8.80%  "
$ws.Range("E51").Value = "  +2.97%  "
